$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Note: "12÷8=1, 4" is both a source and a target text in this batch.
# Replace it first (as a source) before any other rule turns a different
# cell's text INTO "12÷8=1, 4", so we don't clobber the freshly-written value.
Replace-Text "12÷8=1, 4" "60÷3=20, 0"

Replace-Text "58÷9=6, 4" "17÷7=2, 3"
Replace-Text "28÷8=3, 4" "87÷6=14, 3"
Replace-Text "59÷4=14, 3" "94÷4=23, 2"
Replace-Text "47÷5=9, 2" "21÷2=10, 1"
Replace-Text "80÷5=16, 0" "84÷3=28, 0"
Replace-Text "75÷9=8, 3" "57÷4=14, 1"
Replace-Text "41÷4=10, 1" "58÷2=29, 0"
Replace-Text "39÷4=9, 3" "65÷2=32, 1"
Replace-Text "89÷6=14, 5" "16÷2=8, 0"
Replace-Text "54÷8=6, 6" "94÷4=23, 2"
Replace-Text "79÷6=13, 1" "96÷8=12, 0"
Replace-Text "44÷5=8, 4" "92÷5=18, 2"
Replace-Text "42÷6=7, 0" "86÷6=14, 2"
Replace-Text "37÷2=18, 1" "42÷8=5, 2"
Replace-Text "78÷5=15, 3" "12÷8=1, 4"
Replace-Text "85÷8=10, 5" "79÷9=8, 7"
Replace-Text "63÷8=7, 7" "37÷7=5, 2"
Replace-Text "49÷3=16, 1" "95÷3=31, 2"
Replace-Text "71÷9=7, 8" "34÷6=5, 4"
Replace-Text "89÷7=12, 5" "45÷8=5, 5"
Replace-Text "86÷7=12, 2" "40÷2=20, 0"
Replace-Text "31÷6=5, 1" "47÷7=6, 5"
Replace-Text "55÷3=18, 1" "66÷6=11, 0"
Replace-Text "31÷2=15, 1" "93÷8=11, 5"

Write-Output "Replacements complete"
